$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet: row 3 corresponds to b.md. Its zh-cn / de-de status moves
# from "Handed back: in sync with en-US" to "Ready for handoff" and the
# "Latest HO Xliff Generate Date" column is refreshed.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 12:10:19"

# ---------------------------------------------------------------------------
# "zh-cn" sheet: row 3 (b.md) got a new handoff package generated for it.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-31 12:09:59"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c93189dfd7913e1646b53cba674e9947a50c569a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b32704cf07c5f81b6ab7906908e7c66f56c4723/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# "de-de" sheet: same story as zh-cn, but for the German handoff package.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-31 12:10:19"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c93189dfd7913e1646b53cba674e9947a50c569a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b32704cf07c5f81b6ab7906908e7c66f56c4723/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
